# The workbook's single sheet is protected; unprotect so we can write to
# the cells, then re-protect once all edits are applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update confidentiality footer: model date moves from 2021-05-11 to 2021-05-12
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) for holdings rows 2-35
$ws.Range("D2").Value = 0.03606320345147849
$ws.Range("E2").Value = -0.01244651886425518

$ws.Range("D3").Value = 0.02041788633897391
$ws.Range("E3").Value = -0.007794232268121593

$ws.Range("D4").Value = 0.01918882237324876
$ws.Range("E4").Value = -0.001616161616161627

$ws.Range("D5").Value = 0.03718612191119957
$ws.Range("E5").Value = -0.006062767475035646

$ws.Range("D6").Value = 0.0342747460382057
$ws.Range("E6").Value = 0.0004003202562048447

$ws.Range("D7").Value = 0.01977030183910479
$ws.Range("E7").Value = -0.004063467492260053

$ws.Range("D8").Value = 0.03734832407799099
$ws.Range("E8").Value = -0.01926333615580034

$ws.Range("D9").Value = 0.02039891174587756
$ws.Range("E9").Value = -0.009091727428211405

$ws.Range("D10").Value = 0.02642926186257004
$ws.Range("E10").Value = -0.02325581395348841

$ws.Range("D11").Value = 0.02387763815384656
$ws.Range("E11").Value = -0.01682242990654204

$ws.Range("D12").Value = 0.05734221497311155
$ws.Range("E12").Value = -0.01757510984443655

$ws.Range("D13").Value = 0.02489956281160318
$ws.Range("E13").Value = -0.009218289085545672

$ws.Range("D14").Value = 0.0275213211050242
$ws.Range("E14").Value = -0.0218844984802431

$ws.Range("D15").Value = 0.03344261831838332
$ws.Range("E15").Value = -0.002220703792278678

$ws.Range("D16").Value = 0.01946910567718551
$ws.Range("E16").Value = 0.008640406607369533

$ws.Range("D17").Value = 0.03126217233536462
$ws.Range("E17").Value = -0.02373959862946651

$ws.Range("D18").Value = 0.04186652154163367
$ws.Range("E18").Value = -0.007638888888888973

$ws.Range("D19").Value = 0.1253789084696998
$ws.Range("E19").Value = -0.006680026720107035

$ws.Range("D20").Value = 0.009087299884033102
$ws.Range("E20").Value = 0.006971340046475705

$ws.Range("D21").Value = 0.01543629054662054
$ws.Range("E21").Value = -0.0242176115467353

$ws.Range("D22").Value = 0.01701113076664727
$ws.Range("E22").Value = -0.0170251809562646

$ws.Range("D23").Value = 0.0156062457729567
$ws.Range("E23").Value = -0.02030637691485571

$ws.Range("D24").Value = 0.02138686576115522
$ws.Range("E24").Value = -0.02471541380371256

$ws.Range("D25").Value = 0.01237123067094034
$ws.Range("E25").Value = -0.01160220994475136

$ws.Range("D26").Value = 0.04230033582734467
$ws.Range("E26").Value = -0.01497639589776967

$ws.Range("D27").Value = 0.02392777800603397
$ws.Range("E27").Value = -0.0001961168856638995

$ws.Range("D28").Value = 0.04554392010043681
$ws.Range("E28").Value = -0.01335877862595425

$ws.Range("D29").Value = 0.05546538798311986
$ws.Range("E29").Value = -0.02469358327325166

$ws.Range("D30").Value = 0.01292088178708635
$ws.Range("E30").Value = -0.02883355176933156

$ws.Range("D31").Value = 0.0206032456704582
$ws.Range("E31").Value = -0.003453568687643904

$ws.Range("D32").Value = 0.01359973355998791
$ws.Range("E32").Value = -0.03469292076887009

$ws.Range("D33").Value = 0.04187891623551113
$ws.Range("E33").Value = -0.003091190108191699

$ws.Range("D34").Value = 0.01672309440316578
$ws.Range("E34").Value = -0.01120908683305932

$ws.Range("E35").Value = -0.01193527409246031

# Restore sheet protection
$ws.Protect()
